$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three dialogue lines to wrap them in green color BBCode-style tags
$ws.Range("B2").Value = " <color=#00CC00>(Among all the suspects, only one person possesses an item capable of causing such a wound.)</color>"
$ws.Range("B3").Value = " <color=#00CC00>(Moreover, this person’s clothing seems different from when we first saw him/her.)</color>"
$ws.Range("B4").Value = " <color=#00CC00>(Observe all the suspects carefully——identify the accomplice!)</color>"

# Increase row heights for rows 2 and 3 to fit the longer text
$ws.Rows.Item(2).RowHeight = 51
$ws.Rows.Item(3).RowHeight = 51

# Move the active selection to B13
$ws.Range("B13").Select()
